$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1946308724832215
$ws.Range("C2").Value = 0.5536912751677853
$ws.Range("J2").Value = 0.0302013422818792
$ws.Range("P2").Value = 0.1208053691275168
$ws.Range("S2").Value = 0.1006711409395973
$ws.Range("B3").Value = 0.01744186046511628
$ws.Range("C3").Value = 0.02325581395348837
$ws.Range("J3").Value = 0.02325581395348837
$ws.Range("P3").Value = 0.6744186046511628
$ws.Range("S3").Value = 0.2616279069767442
$ws.Range("J4").Value = 0.08571428571428572
$ws.Range("P4").Value = 0.5142857142857142
$ws.Range("S4").Value = 0.4
$ws.Range("B6").Value = 0.07798165137614679
$ws.Range("D6").Value = 0.01376146788990826
$ws.Range("F6").Value = 0.03211009174311927
$ws.Range("J6").Value = 0.2431192660550459
$ws.Range("O6").Value = 0.01834862385321101
$ws.Range("Q6").Value = 0.2064220183486239
$ws.Range("R6").Value = 0.05963302752293578
$ws.Range("S6").Value = 0.3486238532110092
$ws.Range("B7").Value = 0.1421800947867299
$ws.Range("D7").Value = 0.009478672985781991
$ws.Range("E7").Value = 0.004739336492890996
$ws.Range("F7").Value = 0.07109004739336493
$ws.Range("J7").Value = 0.1421800947867299
$ws.Range("O7").Value = 0.01895734597156398
$ws.Range("Q7").Value = 0.1753554502369668
$ws.Range("R7").Value = 0.07582938388625593
$ws.Range("S7").Value = 0.3601895734597156
$ws.Range("B8").Value = 0.0951276102088167
$ws.Range("D8").Value = 0.009280742459396751
$ws.Range("F8").Value = 0.06496519721577726
$ws.Range("J8").Value = 0.1276102088167053
$ws.Range("O8").Value = 0.02320185614849188
$ws.Range("Q8").Value = 0.1531322505800464
$ws.Range("R8").Value = 0.08120649651972157
$ws.Range("S8").Value = 0.4454756380510441
$ws.Range("B9").Value = 0.05813953488372093
$ws.Range("D9").Value = 0.01162790697674419
$ws.Range("F9").Value = 0.0872093023255814
$ws.Range("J9").Value = 0.08139534883720931
$ws.Range("O9").Value = 0.01744186046511628
$ws.Range("Q9").Value = 0.2325581395348837
$ws.Range("R9").Value = 0.0755813953488372
$ws.Range("S9").Value = 0.436046511627907
$ws.Range("B10").Value = 0.1055900621118012
$ws.Range("D10").Value = 0.02018633540372671
$ws.Range("E10").Value = 0.0007763975155279503
$ws.Range("F10").Value = 0.06909937888198758
$ws.Range("J10").Value = 0.1234472049689441
$ws.Range("O10").Value = 0.01475155279503106
$ws.Range("Q10").Value = 0.1940993788819876
$ws.Range("R10").Value = 0.09006211180124224
$ws.Range("S10").Value = 0.3819875776397516
$ws.Range("G11").Value = 0.163323782234957
$ws.Range("J11").Value = 0.0830945558739255
$ws.Range("K11").Value = 0.2063037249283668
$ws.Range("L11").Value = 0.5358166189111748
$ws.Range("S11").Value = 0.01146131805157593
$ws.Range("G12").Value = 0.6770833333333334
$ws.Range("J12").Value = 0.2447916666666667
$ws.Range("K12").Value = 0.005208333333333333
$ws.Range("L12").Value = 0.03645833333333334
$ws.Range("S12").Value = 0.03645833333333334
$ws.Range("G13").Value = 0.7714285714285715
$ws.Range("J13").Value = 0.1142857142857143
$ws.Range("S13").Value = 0.1142857142857143
$ws.Range("F15").Value = 0.0374331550802139
$ws.Range("H15").Value = 0.1657754010695187
$ws.Range("I15").Value = 0.06417112299465241
$ws.Range("J15").Value = 0.3422459893048128
$ws.Range("K15").Value = 0.09625668449197861
$ws.Range("M15").Value = 0.0106951871657754
$ws.Range("N15").Value = 0.0053475935828877
$ws.Range("O15").Value = 0.03208556149732621
$ws.Range("S15").Value = 0.2459893048128342
$ws.Range("F16").Value = 0.01796407185628742
$ws.Range("H16").Value = 0.1736526946107785
$ws.Range("I16").Value = 0.05389221556886228
$ws.Range("J16").Value = 0.4011976047904192
$ws.Range("K16").Value = 0.1377245508982036
$ws.Range("O16").Value = 0.03592814371257485
$ws.Range("S16").Value = 0.1796407185628743
$ws.Range("F17").Value = 0.02097902097902098
$ws.Range("H17").Value = 0.1678321678321678
$ws.Range("I17").Value = 0.09557109557109557
$ws.Range("J17").Value = 0.4242424242424243
$ws.Range("K17").Value = 0.1258741258741259
$ws.Range("M17").Value = 0.01398601398601399
$ws.Range("N17").Value = 0.002331002331002331
$ws.Range("O17").Value = 0.04662004662004662
$ws.Range("S17").Value = 0.1025641025641026
$ws.Range("F18").Value = 0.0155440414507772
$ws.Range("H18").Value = 0.1813471502590674
$ws.Range("I18").Value = 0.1191709844559585
$ws.Range("J18").Value = 0.461139896373057
$ws.Range("K18").Value = 0.06217616580310881
$ws.Range("M18").Value = 0.005181347150259068
$ws.Range("O18").Value = 0.05699481865284974
$ws.Range("S18").Value = 0.09844559585492228
$ws.Range("F19").Value = 0.01226993865030675
$ws.Range("H19").Value = 0.2032208588957055
$ws.Range("I19").Value = 0.06671779141104295
$ws.Range("J19").Value = 0.3826687116564417
$ws.Range("K19").Value = 0.1273006134969325
$ws.Range("M19").Value = 0.01993865030674847
$ws.Range("O19").Value = 0.06058282208588957
$ws.Range("S19").Value = 0.1273006134969325

Write-Output "Updated 109 cells in the team-specific transition matrix."
